$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 1014-1016: new weekly price data (week of 2023-03-28)
$ws.Range("D1014").Value = 45013
$ws.Range("L1014").Value = 'Pintón'
$ws.Range("M1014").Value = 120
$ws.Range("N1014").Value = 21000
$ws.Range("O1014").Value = 21000
$ws.Range("P1014").Value = 21000
$ws.Range("S1014").Value = 1050
$ws.Range("D1015").Value = 45013
$ws.Range("L1015").Value = 'Primera Maduro'
$ws.Range("M1015").Value = 150
$ws.Range("N1015").Value = 24000
$ws.Range("O1015").Value = 24000
$ws.Range("P1015").Value = 24000
$ws.Range("S1015").Value = 1200
$ws.Range("D1016").Value = 45013
$ws.Range("L1016").Value = 'Primera Pintón'
$ws.Range("M1016").Value = 120
$ws.Range("N1016").Value = 25000
$ws.Range("O1016").Value = 25000
$ws.Range("P1016").Value = 25000
$ws.Range("S1016").Value = 1250

# Rows 1017-1088: shift each row block down (data from the row 3 rows below, before-state)
$ws.Range("D1017").Value = 44610
$ws.Range("L1017").Value = 'Pintón'
$ws.Range("M1017").Value = 80
$ws.Range("N1017").Value = 15000
$ws.Range("O1017").Value = 15000
$ws.Range("P1017").Value = 15000
$ws.Range("S1017").Value = 750
$ws.Range("D1018").Value = 44610
$ws.Range("L1018").Value = 'Primera Maduro'
$ws.Range("M1018").Value = 120
$ws.Range("N1018").Value = 17000
$ws.Range("O1018").Value = 17000
$ws.Range("P1018").Value = 17000
$ws.Range("S1018").Value = 850
$ws.Range("D1019").Value = 44610
$ws.Range("L1019").Value = 'Primera Pintón'
$ws.Range("M1019").Value = 120
$ws.Range("N1019").Value = 18000
$ws.Range("O1019").Value = 18000
$ws.Range("P1019").Value = 18000
$ws.Range("S1019").Value = 900
$ws.Range("D1020").Value = 44939
$ws.Range("L1020").Value = 'Pintón'
$ws.Range("M1020").Value = 80
$ws.Range("N1020").Value = 21000
$ws.Range("O1020").Value = 21000
$ws.Range("P1020").Value = 21000
$ws.Range("S1020").Value = 1050
$ws.Range("D1021").Value = 44939
$ws.Range("L1021").Value = 'Primera Maduro'
$ws.Range("M1021").Value = 120
$ws.Range("N1021").Value = 24000
$ws.Range("O1021").Value = 24000
$ws.Range("P1021").Value = 24000
$ws.Range("S1021").Value = 1200
$ws.Range("D1022").Value = 44939
$ws.Range("L1022").Value = 'Primera Pintón'
$ws.Range("M1022").Value = 120
$ws.Range("N1022").Value = 25000
$ws.Range("O1022").Value = 25000
$ws.Range("P1022").Value = 25000
$ws.Range("S1022").Value = 1250
$ws.Range("D1023").Value = 44673
$ws.Range("L1023").Value = 'Pintón'
$ws.Range("M1023").Value = 80
$ws.Range("N1023").Value = 14000
$ws.Range("O1023").Value = 14000
$ws.Range("P1023").Value = 14000
$ws.Range("S1023").Value = 700
$ws.Range("D1024").Value = 44673
$ws.Range("L1024").Value = 'Primera Maduro'
$ws.Range("M1024").Value = 120
$ws.Range("N1024").Value = 16000
$ws.Range("O1024").Value = 16000
$ws.Range("P1024").Value = 16000
$ws.Range("S1024").Value = 800
$ws.Range("D1025").Value = 44673
$ws.Range("L1025").Value = 'Primera Pintón'
$ws.Range("M1025").Value = 120
$ws.Range("N1025").Value = 17000
$ws.Range("O1025").Value = 17000
$ws.Range("P1025").Value = 17000
$ws.Range("S1025").Value = 850
$ws.Range("D1026").Value = 44568
$ws.Range("L1026").Value = 'Pintón'
$ws.Range("M1026").Value = 80
$ws.Range("N1026").Value = 12000
$ws.Range("O1026").Value = 12000
$ws.Range("P1026").Value = 12000
$ws.Range("S1026").Value = 600
$ws.Range("D1027").Value = 44568
$ws.Range("L1027").Value = 'Primera Maduro'
$ws.Range("M1027").Value = 120
$ws.Range("N1027").Value = 14000
$ws.Range("O1027").Value = 14000
$ws.Range("P1027").Value = 14000
$ws.Range("S1027").Value = 700
$ws.Range("D1028").Value = 44568
$ws.Range("L1028").Value = 'Primera Pintón'
$ws.Range("M1028").Value = 120
$ws.Range("N1028").Value = 15000
$ws.Range("O1028").Value = 15000
$ws.Range("P1028").Value = 15000
$ws.Range("S1028").Value = 750
$ws.Range("D1029").Value = 44473
$ws.Range("L1029").Value = 'Pintón'
$ws.Range("M1029").Value = 80
$ws.Range("N1029").Value = 21000
$ws.Range("O1029").Value = 21000
$ws.Range("P1029").Value = 21000
$ws.Range("S1029").Value = 1050
$ws.Range("D1030").Value = 44473
$ws.Range("L1030").Value = 'Primera Maduro'
$ws.Range("M1030").Value = 120
$ws.Range("N1030").Value = 22000
$ws.Range("O1030").Value = 22000
$ws.Range("P1030").Value = 22000
$ws.Range("S1030").Value = 1100
$ws.Range("D1031").Value = 44473
$ws.Range("L1031").Value = 'Primera Pintón'
$ws.Range("M1031").Value = 120
$ws.Range("N1031").Value = 23000
$ws.Range("O1031").Value = 23000
$ws.Range("P1031").Value = 23000
$ws.Range("S1031").Value = 1150
$ws.Range("D1032").Value = 44620
$ws.Range("L1032").Value = 'Pintón'
$ws.Range("M1032").Value = 80
$ws.Range("N1032").Value = 16000
$ws.Range("O1032").Value = 16000
$ws.Range("P1032").Value = 16000
$ws.Range("S1032").Value = 800
$ws.Range("D1033").Value = 44620
$ws.Range("L1033").Value = 'Primera Maduro'
$ws.Range("M1033").Value = 120
$ws.Range("N1033").Value = 18000
$ws.Range("O1033").Value = 18000
$ws.Range("P1033").Value = 18000
$ws.Range("S1033").Value = 900
$ws.Range("D1034").Value = 44620
$ws.Range("L1034").Value = 'Primera Pintón'
$ws.Range("M1034").Value = 120
$ws.Range("N1034").Value = 19000
$ws.Range("O1034").Value = 19000
$ws.Range("P1034").Value = 19000
$ws.Range("S1034").Value = 950
$ws.Range("D1035").Value = 44413
$ws.Range("L1035").Value = 'Pintón'
$ws.Range("M1035").Value = 80
$ws.Range("N1035").Value = 14000
$ws.Range("O1035").Value = 14000
$ws.Range("P1035").Value = 14000
$ws.Range("S1035").Value = 700
$ws.Range("D1036").Value = 44413
$ws.Range("L1036").Value = 'Primera Maduro'
$ws.Range("M1036").Value = 120
$ws.Range("N1036").Value = 15500
$ws.Range("O1036").Value = 15500
$ws.Range("P1036").Value = 15500
$ws.Range("S1036").Value = 775
$ws.Range("D1037").Value = 44413
$ws.Range("L1037").Value = 'Primera Pintón'
$ws.Range("M1037").Value = 120
$ws.Range("N1037").Value = 16000
$ws.Range("O1037").Value = 16000
$ws.Range("P1037").Value = 16000
$ws.Range("S1037").Value = 800
$ws.Range("D1038").Value = 44606
$ws.Range("L1038").Value = 'Pintón'
$ws.Range("M1038").Value = 80
$ws.Range("N1038").Value = 15000
$ws.Range("O1038").Value = 15000
$ws.Range("P1038").Value = 15000
$ws.Range("S1038").Value = 750
$ws.Range("D1039").Value = 44606
$ws.Range("L1039").Value = 'Primera Maduro'
$ws.Range("M1039").Value = 120
$ws.Range("N1039").Value = 17000
$ws.Range("O1039").Value = 17000
$ws.Range("P1039").Value = 17000
$ws.Range("S1039").Value = 850
$ws.Range("D1040").Value = 44606
$ws.Range("L1040").Value = 'Primera Pintón'
$ws.Range("M1040").Value = 120
$ws.Range("N1040").Value = 18000
$ws.Range("O1040").Value = 18000
$ws.Range("P1040").Value = 18000
$ws.Range("S1040").Value = 900
$ws.Range("D1041").Value = 44567
$ws.Range("L1041").Value = 'Pintón'
$ws.Range("M1041").Value = 80
$ws.Range("N1041").Value = 12000
$ws.Range("O1041").Value = 12000
$ws.Range("P1041").Value = 12000
$ws.Range("S1041").Value = 600
$ws.Range("D1042").Value = 44567
$ws.Range("L1042").Value = 'Primera Maduro'
$ws.Range("M1042").Value = 120
$ws.Range("N1042").Value = 14000
$ws.Range("O1042").Value = 14000
$ws.Range("P1042").Value = 14000
$ws.Range("S1042").Value = 700
$ws.Range("D1043").Value = 44567
$ws.Range("L1043").Value = 'Primera Pintón'
$ws.Range("M1043").Value = 120
$ws.Range("N1043").Value = 15000
$ws.Range("O1043").Value = 15000
$ws.Range("P1043").Value = 15000
$ws.Range("S1043").Value = 750
$ws.Range("D1044").Value = 44960
$ws.Range("L1044").Value = 'Maduro'
$ws.Range("M1044").Value = 80
$ws.Range("N1044").Value = 19000
$ws.Range("O1044").Value = 19000
$ws.Range("P1044").Value = 19000
$ws.Range("S1044").Value = 950
$ws.Range("D1045").Value = 44960
$ws.Range("L1045").Value = 'Primera Maduro'
$ws.Range("M1045").Value = 120
$ws.Range("N1045").Value = 22000
$ws.Range("O1045").Value = 22000
$ws.Range("P1045").Value = 22000
$ws.Range("S1045").Value = 1100
$ws.Range("D1046").Value = 44960
$ws.Range("L1046").Value = 'Primera Pintón'
$ws.Range("M1046").Value = 120
$ws.Range("N1046").Value = 23000
$ws.Range("O1046").Value = 23000
$ws.Range("P1046").Value = 23000
$ws.Range("S1046").Value = 1150
$ws.Range("D1047").Value = 44364
$ws.Range("L1047").Value = 'Pintón'
$ws.Range("M1047").Value = 80
$ws.Range("N1047").Value = 18000
$ws.Range("O1047").Value = 18000
$ws.Range("P1047").Value = 18000
$ws.Range("S1047").Value = 900
$ws.Range("D1048").Value = 44364
$ws.Range("L1048").Value = 'Primera Maduro'
$ws.Range("M1048").Value = 120
$ws.Range("N1048").Value = 19500
$ws.Range("O1048").Value = 19500
$ws.Range("P1048").Value = 19500
$ws.Range("S1048").Value = 975
$ws.Range("D1049").Value = 44364
$ws.Range("L1049").Value = 'Primera Pintón'
$ws.Range("M1049").Value = 120
$ws.Range("N1049").Value = 20000
$ws.Range("O1049").Value = 20000
$ws.Range("P1049").Value = 20000
$ws.Range("S1049").Value = 1000
$ws.Range("D1050").Value = 44168
$ws.Range("L1050").Value = 'Pintón'
$ws.Range("M1050").Value = 80
$ws.Range("N1050").Value = 19000
$ws.Range("O1050").Value = 19000
$ws.Range("P1050").Value = 19000
$ws.Range("S1050").Value = 950
$ws.Range("D1051").Value = 44168
$ws.Range("L1051").Value = 'Primera Maduro'
$ws.Range("M1051").Value = 120
$ws.Range("N1051").Value = 20500
$ws.Range("O1051").Value = 20500
$ws.Range("P1051").Value = 20500
$ws.Range("S1051").Value = 1025
$ws.Range("D1052").Value = 44168
$ws.Range("L1052").Value = 'Primera Pintón'
$ws.Range("M1052").Value = 120
$ws.Range("N1052").Value = 21000
$ws.Range("O1052").Value = 21000
$ws.Range("P1052").Value = 21000
$ws.Range("S1052").Value = 1050
$ws.Range("D1053").Value = 44677
$ws.Range("L1053").Value = 'Pintón'
$ws.Range("M1053").Value = 120
$ws.Range("N1053").Value = 14000
$ws.Range("O1053").Value = 14000
$ws.Range("P1053").Value = 14000
$ws.Range("S1053").Value = 700
$ws.Range("D1054").Value = 44677
$ws.Range("L1054").Value = 'Primera Maduro'
$ws.Range("M1054").Value = 160
$ws.Range("N1054").Value = 16000
$ws.Range("O1054").Value = 16000
$ws.Range("P1054").Value = 16000
$ws.Range("S1054").Value = 800
$ws.Range("D1055").Value = 44677
$ws.Range("L1055").Value = 'Primera Pintón'
$ws.Range("M1055").Value = 160
$ws.Range("N1055").Value = 17000
$ws.Range("O1055").Value = 17000
$ws.Range("P1055").Value = 17000
$ws.Range("S1055").Value = 850
$ws.Range("D1056").Value = 44747
$ws.Range("L1056").Value = 'Pintón'
$ws.Range("M1056").Value = 120
$ws.Range("N1056").Value = 23000
$ws.Range("O1056").Value = 23000
$ws.Range("P1056").Value = 23000
$ws.Range("S1056").Value = 1150
$ws.Range("D1057").Value = 44747
$ws.Range("L1057").Value = 'Primera Maduro'
$ws.Range("M1057").Value = 120
$ws.Range("N1057").Value = 25000
$ws.Range("O1057").Value = 25000
$ws.Range("P1057").Value = 25000
$ws.Range("S1057").Value = 1250
$ws.Range("D1058").Value = 44747
$ws.Range("L1058").Value = 'Primera Pintón'
$ws.Range("M1058").Value = 160
$ws.Range("N1058").Value = 26000
$ws.Range("O1058").Value = 26000
$ws.Range("P1058").Value = 26000
$ws.Range("S1058").Value = 1300
$ws.Range("D1059").Value = 44771
$ws.Range("L1059").Value = 'Pintón'
$ws.Range("M1059").Value = 80
$ws.Range("N1059").Value = 30000
$ws.Range("O1059").Value = 30000
$ws.Range("P1059").Value = 30000
$ws.Range("S1059").Value = 1500
$ws.Range("D1060").Value = 44771
$ws.Range("L1060").Value = 'Primera Maduro'
$ws.Range("M1060").Value = 120
$ws.Range("N1060").Value = 32000
$ws.Range("O1060").Value = 32000
$ws.Range("P1060").Value = 32000
$ws.Range("S1060").Value = 1600
$ws.Range("D1061").Value = 44771
$ws.Range("L1061").Value = 'Primera Pintón'
$ws.Range("M1061").Value = 120
$ws.Range("N1061").Value = 33000
$ws.Range("O1061").Value = 33000
$ws.Range("P1061").Value = 33000
$ws.Range("S1061").Value = 1650
$ws.Range("D1062").Value = 44245
$ws.Range("L1062").Value = 'Pintón'
$ws.Range("M1062").Value = 80
$ws.Range("N1062").Value = 13000
$ws.Range("O1062").Value = 13000
$ws.Range("P1062").Value = 13000
$ws.Range("S1062").Value = 650
$ws.Range("D1063").Value = 44245
$ws.Range("L1063").Value = 'Primera Maduro'
$ws.Range("M1063").Value = 120
$ws.Range("N1063").Value = 14500
$ws.Range("O1063").Value = 14500
$ws.Range("P1063").Value = 14500
$ws.Range("S1063").Value = 725
$ws.Range("D1064").Value = 44245
$ws.Range("L1064").Value = 'Primera Pintón'
$ws.Range("M1064").Value = 120
$ws.Range("N1064").Value = 15000
$ws.Range("O1064").Value = 15000
$ws.Range("P1064").Value = 15000
$ws.Range("S1064").Value = 750
$ws.Range("D1065").Value = 44498
$ws.Range("L1065").Value = 'Pintón'
$ws.Range("M1065").Value = 80
$ws.Range("N1065").Value = 21000
$ws.Range("O1065").Value = 21000
$ws.Range("P1065").Value = 21000
$ws.Range("S1065").Value = 1050
$ws.Range("D1066").Value = 44498
$ws.Range("L1066").Value = 'Primera Maduro'
$ws.Range("M1066").Value = 120
$ws.Range("N1066").Value = 22500
$ws.Range("O1066").Value = 22500
$ws.Range("P1066").Value = 22500
$ws.Range("S1066").Value = 1125
$ws.Range("D1067").Value = 44498
$ws.Range("L1067").Value = 'Primera Pintón'
$ws.Range("M1067").Value = 120
$ws.Range("N1067").Value = 23000
$ws.Range("O1067").Value = 23000
$ws.Range("P1067").Value = 23000
$ws.Range("S1067").Value = 1150
$ws.Range("D1068").Value = 44536
$ws.Range("L1068").Value = 'Pintón'
$ws.Range("M1068").Value = 80
$ws.Range("N1068").Value = 16000
$ws.Range("O1068").Value = 16000
$ws.Range("P1068").Value = 16000
$ws.Range("S1068").Value = 800
$ws.Range("D1069").Value = 44536
$ws.Range("L1069").Value = 'Primera Maduro'
$ws.Range("M1069").Value = 120
$ws.Range("N1069").Value = 18000
$ws.Range("O1069").Value = 18000
$ws.Range("P1069").Value = 18000
$ws.Range("S1069").Value = 900
$ws.Range("D1070").Value = 44536
$ws.Range("L1070").Value = 'Primera Pintón'
$ws.Range("M1070").Value = 120
$ws.Range("N1070").Value = 19000
$ws.Range("O1070").Value = 19000
$ws.Range("P1070").Value = 19000
$ws.Range("S1070").Value = 950
$ws.Range("D1071").Value = 44242
$ws.Range("L1071").Value = 'Pintón'
$ws.Range("M1071").Value = 80
$ws.Range("N1071").Value = 13000
$ws.Range("O1071").Value = 13000
$ws.Range("P1071").Value = 13000
$ws.Range("S1071").Value = 650
$ws.Range("D1072").Value = 44242
$ws.Range("L1072").Value = 'Primera Maduro'
$ws.Range("M1072").Value = 120
$ws.Range("N1072").Value = 14500
$ws.Range("O1072").Value = 14500
$ws.Range("P1072").Value = 14500
$ws.Range("S1072").Value = 725
$ws.Range("D1073").Value = 44242
$ws.Range("L1073").Value = 'Primera Pintón'
$ws.Range("M1073").Value = 120
$ws.Range("N1073").Value = 15000
$ws.Range("O1073").Value = 15000
$ws.Range("P1073").Value = 15000
$ws.Range("S1073").Value = 750
$ws.Range("D1074").Value = 44988
$ws.Range("L1074").Value = 'Pintón'
$ws.Range("M1074").Value = 80
$ws.Range("N1074").Value = 22000
$ws.Range("O1074").Value = 22000
$ws.Range("P1074").Value = 22000
$ws.Range("S1074").Value = 1100
$ws.Range("D1075").Value = 44988
$ws.Range("L1075").Value = 'Primera Maduro'
$ws.Range("M1075").Value = 120
$ws.Range("N1075").Value = 24000
$ws.Range("O1075").Value = 24000
$ws.Range("P1075").Value = 24000
$ws.Range("S1075").Value = 1200
$ws.Range("D1076").Value = 44988
$ws.Range("L1076").Value = 'Primera Pintón'
$ws.Range("M1076").Value = 120
$ws.Range("N1076").Value = 25000
$ws.Range("O1076").Value = 25000
$ws.Range("P1076").Value = 25000
$ws.Range("S1076").Value = 1250
$ws.Range("D1077").Value = 44608
$ws.Range("L1077").Value = 'Pintón'
$ws.Range("M1077").Value = 120
$ws.Range("N1077").Value = 15000
$ws.Range("O1077").Value = 15000
$ws.Range("P1077").Value = 15000
$ws.Range("S1077").Value = 750
$ws.Range("D1078").Value = 44608
$ws.Range("L1078").Value = 'Primera Maduro'
$ws.Range("M1078").Value = 160
$ws.Range("N1078").Value = 17000
$ws.Range("O1078").Value = 17000
$ws.Range("P1078").Value = 17000
$ws.Range("S1078").Value = 850
$ws.Range("D1079").Value = 44608
$ws.Range("L1079").Value = 'Primera Pintón'
$ws.Range("M1079").Value = 160
$ws.Range("N1079").Value = 18000
$ws.Range("O1079").Value = 18000
$ws.Range("P1079").Value = 18000
$ws.Range("S1079").Value = 900
$ws.Range("D1080").Value = 44414
$ws.Range("L1080").Value = 'Pintón'
$ws.Range("M1080").Value = 80
$ws.Range("N1080").Value = 13500
$ws.Range("O1080").Value = 13500
$ws.Range("P1080").Value = 13500
$ws.Range("S1080").Value = 675
$ws.Range("D1081").Value = 44414
$ws.Range("L1081").Value = 'Primera Maduro'
$ws.Range("M1081").Value = 120
$ws.Range("N1081").Value = 15000
$ws.Range("O1081").Value = 15000
$ws.Range("P1081").Value = 15000
$ws.Range("S1081").Value = 750
$ws.Range("D1082").Value = 44414
$ws.Range("L1082").Value = 'Primera Pintón'
$ws.Range("M1082").Value = 120
$ws.Range("N1082").Value = 15500
$ws.Range("O1082").Value = 15500
$ws.Range("P1082").Value = 15500
$ws.Range("S1082").Value = 775
$ws.Range("D1083").Value = 44925
$ws.Range("L1083").Value = 'Pintón'
$ws.Range("M1083").Value = 120
$ws.Range("N1083").Value = 18000
$ws.Range("O1083").Value = 18000
$ws.Range("P1083").Value = 18000
$ws.Range("S1083").Value = 900
$ws.Range("D1084").Value = 44925
$ws.Range("L1084").Value = 'Primera Maduro'
$ws.Range("M1084").Value = 120
$ws.Range("N1084").Value = 20000
$ws.Range("O1084").Value = 20000
$ws.Range("P1084").Value = 20000
$ws.Range("S1084").Value = 1000
$ws.Range("D1085").Value = 44925
$ws.Range("L1085").Value = 'Primera Pintón'
$ws.Range("M1085").Value = 160
$ws.Range("N1085").Value = 21000
$ws.Range("O1085").Value = 21000
$ws.Range("P1085").Value = 21000
$ws.Range("S1085").Value = 1050
$ws.Range("D1086").Value = 44210
$ws.Range("L1086").Value = 'Pintón'
$ws.Range("M1086").Value = 80
$ws.Range("N1086").Value = 14000
$ws.Range("O1086").Value = 14000
$ws.Range("P1086").Value = 14000
$ws.Range("S1086").Value = 700
$ws.Range("D1087").Value = 44210
$ws.Range("L1087").Value = 'Primera Maduro'
$ws.Range("M1087").Value = 120
$ws.Range("N1087").Value = 15500
$ws.Range("O1087").Value = 15500
$ws.Range("P1087").Value = 15500
$ws.Range("S1087").Value = 775
$ws.Range("D1088").Value = 44210
$ws.Range("L1088").Value = 'Primera Pintón'
$ws.Range("M1088").Value = 120
$ws.Range("N1088").Value = 16000
$ws.Range("O1088").Value = 16000
$ws.Range("P1088").Value = 16000
$ws.Range("S1088").Value = 800

# New rows 1089-1091: appended, carrying the old last block (rows 1086-1088) down
$ws.Range("A1089").Value = 8
$ws.Range("B1089").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C1089").Value = 'Coquimbo'
$ws.Range("D1089").Value = 44601
$ws.Range("E1089").Value = 4
$ws.Range("F1089").Value = 'Fruta'
$ws.Range("G1089").Value = 100108
$ws.Range("H1089").Value = 'Tropicales y subtropicales'
$ws.Range("I1089").Value = 100108006
$ws.Range("J1089").Value = 'Plátano'
$ws.Range("K1089").Value = 'Sin especificar'
$ws.Range("L1089").Value = 'Pintón'
$ws.Range("M1089").Value = 120
$ws.Range("N1089").Value = 13000
$ws.Range("O1089").Value = 13000
$ws.Range("P1089").Value = 13000
$ws.Range("Q1089").Value = '$/caja 20 kilos'
$ws.Range("R1089").Value = 'Ecuador'
$ws.Range("S1089").Value = 650
$ws.Range("T1089").Value = 20
$ws.Range("D1089").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1090").Value = 8
$ws.Range("B1090").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C1090").Value = 'Coquimbo'
$ws.Range("D1090").Value = 44601
$ws.Range("E1090").Value = 4
$ws.Range("F1090").Value = 'Fruta'
$ws.Range("G1090").Value = 100108
$ws.Range("H1090").Value = 'Tropicales y subtropicales'
$ws.Range("I1090").Value = 100108006
$ws.Range("J1090").Value = 'Plátano'
$ws.Range("K1090").Value = 'Sin especificar'
$ws.Range("L1090").Value = 'Primera Maduro'
$ws.Range("M1090").Value = 160
$ws.Range("N1090").Value = 15000
$ws.Range("O1090").Value = 15000
$ws.Range("P1090").Value = 15000
$ws.Range("Q1090").Value = '$/caja 20 kilos'
$ws.Range("R1090").Value = 'Ecuador'
$ws.Range("S1090").Value = 750
$ws.Range("T1090").Value = 20
$ws.Range("D1090").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A1091").Value = 8
$ws.Range("B1091").Value = 'Terminal La Palmera de La Serena'
$ws.Range("C1091").Value = 'Coquimbo'
$ws.Range("D1091").Value = 44601
$ws.Range("E1091").Value = 4
$ws.Range("F1091").Value = 'Fruta'
$ws.Range("G1091").Value = 100108
$ws.Range("H1091").Value = 'Tropicales y subtropicales'
$ws.Range("I1091").Value = 100108006
$ws.Range("J1091").Value = 'Plátano'
$ws.Range("K1091").Value = 'Sin especificar'
$ws.Range("L1091").Value = 'Primera Pintón'
$ws.Range("M1091").Value = 160
$ws.Range("N1091").Value = 16000
$ws.Range("O1091").Value = 16000
$ws.Range("P1091").Value = 16000
$ws.Range("Q1091").Value = '$/caja 20 kilos'
$ws.Range("R1091").Value = 'Ecuador'
$ws.Range("S1091").Value = 800
$ws.Range("T1091").Value = 20
$ws.Range("D1091").NumberFormat = "YYYY-MM-DD HH:MM:SS"

